$d = $word.ActiveDocument

# Fix typo: PEMILIHIAN -> PEMILIHAN
$d.Content.Find.Execute("PEMILIHIAN", $true, $false, $false, $false, $false, $true, 1, $false, "PEMILIHAN", 2)
